# Auto-generated edit script applying the scheduled-runner cell updates
# described in the commit diff for Sheets/Halicarnassus_Profits.xlsx.
# Each FFXIV leve-profit sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) holds
# static market-price snapshots in columns H:N; this refreshes them.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 262.25
$ws.Range("I2").Value = 173.75
$ws.Range("J2").Value = 439.25
$ws.Range("K2").Value = 173.75
$ws.Range("L2").Value = 439.25
$ws.Range("M2").Value = -60.75
$ws.Range("N2").Value = -665.25
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").ClearContents()
$ws.Range("N3").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H40").Value = 4783.0713
$ws.Range("I40").Value = 3949.8333
$ws.Range("J40").Value = 6866.1665
$ws.Range("K40").Value = 3949.8333
$ws.Range("L40").Value = 6866.1665
$ws.Range("M40").Value = -3774.8333
$ws.Range("N40").Value = -7216.1665
$ws.Range("H86").Value = 4999.5
$ws.Range("I86").Value = 999
$ws.Range("K86").Value = 999
$ws.Range("M86").Value = 124
$ws.Range("H89").Value = 4999.5
$ws.Range("I89").Value = 999
$ws.Range("K89").Value = 4995
$ws.Range("M89").Value = 621
$ws.Range("H96").Value = 430.13333
$ws.Range("I96").Value = 158.8
$ws.Range("J96").Value = 972.8
$ws.Range("K96").Value = 476.4
$ws.Range("L96").Value = 2918.4
$ws.Range("M96").Value = 896.5999999999999
$ws.Range("N96").Value = -5664.4
$ws.Range("H99").Value = 3726.6667
$ws.Range("I99").Value = 1878.6666
$ws.Range("K99").Value = 5635.9998
$ws.Range("M99").Value = -4137.9998
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").ClearContents()
$ws.Range("N102").Value = 0
$ws.Range("H112").Value = 1809.5
$ws.Range("J112").Value = 1979.7273
$ws.Range("L112").Value = 5939.1819
$ws.Range("N112").Value = -8155.1819
$ws.Range("H116").Value = 3500
$ws.Range("I116").Value = 3000
$ws.Range("J116").Value = 3750
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 3750
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -10634
$ws.Range("H138").Value = 4648.3335
$ws.Range("I138").Value = 1900
$ws.Range("J138").Value = 4898.1816
$ws.Range("K138").Value = 5700
$ws.Range("L138").Value = 14694.5448
$ws.Range("M138").Value = -560
$ws.Range("N138").Value = -24974.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 951.1667
$ws.Range("I2").Value = 851.4
$ws.Range("K2").Value = 851.4
$ws.Range("M2").Value = -738.4
$ws.Range("H24").Value = 5025000
$ws.Range("J24").Value = 5025000
$ws.Range("L24").Value = 5025000
$ws.Range("N24").Value = -5025748
$ws.Range("H92").Value = 29846.334
$ws.Range("J92").Value = 29846.334
$ws.Range("L92").Value = 29846.334
$ws.Range("N92").Value = -34838.334
$ws.Range("H96").Value = 2871557.2
$ws.Range("J96").Value = 2871557.2
$ws.Range("L96").Value = 2871557.2
$ws.Range("N96").Value = -2877049.2
$ws.Range("H100").Value = 5025000
$ws.Range("J100").Value = 5025000
$ws.Range("L100").Value = 5025000
$ws.Range("N100").Value = -5027164
$ws.Range("H102").Value = 8330
$ws.Range("I102").Value = 7077.5
$ws.Range("J102").Value = 10000
$ws.Range("K102").Value = 7077.5
$ws.Range("L102").Value = 10000
$ws.Range("M102").Value = -5455.5
$ws.Range("N102").Value = -13244
$ws.Range("H116").Value = 951.1667
$ws.Range("I116").Value = 851.4
$ws.Range("K116").Value = 851.4
$ws.Range("M116").Value = 1442.6
$ws.Range("H140").Value = 100429
$ws.Range("J140").Value = 100429
$ws.Range("L140").Value = 100429
$ws.Range("N140").Value = -110789

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 951.1667
$ws.Range("I3").Value = 851.4
$ws.Range("K3").Value = 851.4
$ws.Range("M3").Value = -737.4
$ws.Range("H86").Value = 5057.7085
$ws.Range("I86").Value = 3928.1333
$ws.Range("K86").Value = 3928.1333
$ws.Range("M86").Value = -2805.1333
$ws.Range("H89").Value = 5057.7085
$ws.Range("I89").Value = 3928.1333
$ws.Range("K89").Value = 19640.6665
$ws.Range("M89").Value = -14024.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 27590
$ws.Range("I69").Value = 20180
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 20180
$ws.Range("L69").Value = 35000
$ws.Range("M69").Value = -19431
$ws.Range("N69").Value = -36498
$ws.Range("H72").Value = 27590
$ws.Range("I72").Value = 20180
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 60540
$ws.Range("L72").Value = 105000
$ws.Range("M72").Value = -56796
$ws.Range("N72").Value = -112488
$ws.Range("H132").Value = 995
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 722.46155
$ws.Range("I34").Value = 82.166664
$ws.Range("J34").Value = 1271.2858
$ws.Range("K34").Value = 246.499992
$ws.Range("L34").Value = 3813.8574
$ws.Range("M34").Value = -162.499992
$ws.Range("N34").Value = -3981.8574
$ws.Range("H39").Value = 4036
$ws.Range("J39").Value = 5099.143
$ws.Range("L39").Value = 15297.429
$ws.Range("N39").Value = -15885.429
$ws.Range("H55").Value = 4801.636
$ws.Range("I55").Value = 1001.5
$ws.Range("J55").Value = 5646.1113
$ws.Range("K55").Value = 3004.5
$ws.Range("L55").Value = 16938.3339
$ws.Range("M55").Value = -2827.5
$ws.Range("N55").Value = -17292.3339
$ws.Range("H68").Value = 1632
$ws.Range("I68").Value = 1632
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 4896
$ws.Range("L68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4085
$ws.Range("H71").Value = 1632
$ws.Range("I71").Value = 1632
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 14688
$ws.Range("L71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -10632
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").ClearContents()
$ws.Range("N100").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 13532235
$ws.Range("I11").Value = 13849230
$ws.Range("K11").Value = 13849230
$ws.Range("M11").Value = -13849091
$ws.Range("H98").Value = 900
$ws.Range("J98").Value = 900
$ws.Range("L98").Value = 900
$ws.Range("N98").Value = -6890
$ws.Range("H102").Value = 1526.0588
$ws.Range("I102").Value = 1309
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 1309
$ws.Range("L102").Value = 4999
$ws.Range("M102").Value = 313
$ws.Range("N102").Value = -8243
$ws.Range("H105").Value = 28166.666
$ws.Range("J105").Value = 28166.666
$ws.Range("L105").Value = 28166.666
$ws.Range("N105").Value = -35154.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9999.5
$ws.Range("I7").Value = 9999.5
$ws.Range("K7").Value = 9999.5
$ws.Range("M7").Value = -9887.5
$ws.Range("H22").Value = 863.3333
$ws.Range("I22").Value = 863.3333
$ws.Range("K22").Value = 863.3333
$ws.Range("M22").Value = -568.3333
$ws.Range("H27").Value = 863.3333
$ws.Range("I27").Value = 863.3333
$ws.Range("K27").Value = 863.3333
$ws.Range("M27").Value = -756.3333
$ws.Range("H104").Value = 12000
$ws.Range("J104").Value = 12000
$ws.Range("L104").Value = 12000
$ws.Range("N104").Value = -18988
$ws.Range("H126").Value = 9999.5
$ws.Range("I126").Value = 9999.5
$ws.Range("K126").Value = 29998.5
$ws.Range("M126").Value = -27528.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 5000
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5346
$ws.Range("H75").Value = 20000
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 20000
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

